# Auto-generated Excel COM-interop script to apply the recorded market-price refresh.
# For each affected leve-profit row, push the new currentAveragePrice / LevePrice /
# LeveProfit figures into columns H-N, matching the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1790695.1
$ws.Range("J17").Value = 1790695.1
$ws.Range("L17").Value = 5372085.300000001
$ws.Range("N17").Value = -5372421.300000001

$ws.Range("H40").Value = 1798.85
$ws.Range("I40").Value = 1651.5
$ws.Range("J40").Value = 2019.875
$ws.Range("K40").Value = 1651.5
$ws.Range("L40").Value = 2019.875
$ws.Range("M40").Value = -1476.5
$ws.Range("N40").Value = -2369.875

$ws.Range("H86").Value = 51479.3
$ws.Range("I86").Value = 21598.6
$ws.Range("J86").Value = 81360
$ws.Range("K86").Value = 21598.6
$ws.Range("L86").Value = 81360
$ws.Range("M86").Value = -20475.6
$ws.Range("N86").Value = -83606

$ws.Range("H89").Value = 51479.3
$ws.Range("I89").Value = 21598.6
$ws.Range("J89").Value = 81360
$ws.Range("K89").Value = 107993
$ws.Range("L89").Value = 406800
$ws.Range("M89").Value = -102377
$ws.Range("N89").Value = -418032

$ws.Range("H100").Value = 4372.778
$ws.Range("I100").Value = 3071
$ws.Range("J100").Value = 6000
$ws.Range("K100").Value = 3071
$ws.Range("L100").Value = 6000
$ws.Range("M100").Value = -2530
$ws.Range("N100").Value = -7082

$ws.Range("H116").Value = 4155.3335
$ws.Range("I116").Value = 3328.25
$ws.Range("J116").Value = 4982.4165
$ws.Range("K116").Value = 3328.25
$ws.Range("L116").Value = 4982.4165
$ws.Range("M116").Value = 113.75
$ws.Range("N116").Value = -11866.4165

$ws.Range("H127").Value = 910934.5600000001
$ws.Range("I127").Value = 1433.3334
$ws.Range("J127").Value = 1251997.5
$ws.Range("K127").Value = 4300.0002
$ws.Range("L127").Value = 3755992.5
$ws.Range("M127").Value = 659.9997999999996
$ws.Range("N127").Value = -3765912.5

$ws.Range("H132").Value = 14954.227
$ws.Range("I132").Value = 15832.3
$ws.Range("J132").Value = 2661.2
$ws.Range("K132").Value = 47496.89999999999
$ws.Range("L132").Value = 7983.599999999999
$ws.Range("M132").Value = -44966.89999999999
$ws.Range("N132").Value = -13043.6

$ws.Range("H138").Value = 6068.519
$ws.Range("I138").Value = 3939.875
$ws.Range("J138").Value = 6609.127
$ws.Range("K138").Value = 11819.625
$ws.Range("L138").Value = 19827.381
$ws.Range("M138").Value = -6679.625
$ws.Range("N138").Value = -30107.381

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1848.625
$ws.Range("I102").Value = 1341.8
$ws.Range("J102").Value = 2693.3333
$ws.Range("K102").Value = 1341.8
$ws.Range("L102").Value = 2693.3333
$ws.Range("M102").Value = 280.2
$ws.Range("N102").Value = -5937.3333

$ws.Range("H139").Value = 62401.11
$ws.Range("J139").Value = 62401.11
$ws.Range("L139").Value = 62401.11
$ws.Range("N139").Value = -72681.11

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 22732
$ws.Range("J52").Value = 22732
$ws.Range("L52").Value = 22732
$ws.Range("N52").Value = -23258

$ws.Range("H86").Value = 3030
$ws.Range("I86").Value = 3450
$ws.Range("J86").Value = 2750
$ws.Range("K86").Value = 3450
$ws.Range("L86").Value = 2750
$ws.Range("M86").Value = -2327
$ws.Range("N86").Value = -4996

$ws.Range("H89").Value = 3030
$ws.Range("I89").Value = 3450
$ws.Range("J89").Value = 2750
$ws.Range("K89").Value = 17250
$ws.Range("L89").Value = 13750
$ws.Range("M89").Value = -11634
$ws.Range("N89").Value = -24982

$ws.Range("H105").Value = 2475
$ws.Range("I105").Value = 2366
$ws.Range("K105").Value = 2366
$ws.Range("M105").Value = -619

$ws.Range("H121").Value = 22732
$ws.Range("J121").Value = 22732
$ws.Range("L121").Value = 22732
$ws.Range("N121").Value = -26226

$ws.Range("H134").Value = 2104.0588
$ws.Range("I134").Value = 1730.75
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 5192.25
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -2657.25
$ws.Range("N134").Value = -14070

$ws.Range("H140").Value = 59830
$ws.Range("J140").Value = 59830
$ws.Range("L140").Value = 59830
$ws.Range("N140").Value = -70190

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6183.2
$ws.Range("I31").Value = 2391
$ws.Range("J31").Value = 14457.091
$ws.Range("K31").Value = 2391
$ws.Range("L31").Value = 14457.091
$ws.Range("M31").Value = -2096
$ws.Range("N31").Value = -15047.091

$ws.Range("H34").Value = 6183.2
$ws.Range("I34").Value = 2391
$ws.Range("J34").Value = 14457.091
$ws.Range("K34").Value = 2391
$ws.Range("L34").Value = 14457.091
$ws.Range("M34").Value = -2189
$ws.Range("N34").Value = -14861.091

$ws.Range("H88").Value = 36150
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 36150
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 36150
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -36962

$ws.Range("H91").Value = 36150
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 36150
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 36150
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -38958

$ws.Range("H132").Value = 3032.2632
$ws.Range("I132").Value = 1752.4
$ws.Range("J132").Value = 4454.3335
$ws.Range("K132").Value = 5257.200000000001
$ws.Range("L132").Value = 13363.0005
$ws.Range("M132").Value = -2727.200000000001
$ws.Range("N132").Value = -18423.0005

$ws.Range("H140").Value = 76382.42999999999
$ws.Range("J140").Value = 76382.42999999999
$ws.Range("L140").Value = 76382.42999999999
$ws.Range("N140").Value = -86742.42999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1681
$ws.Range("I13").Value = 935.8570999999999
$ws.Range("J13").Value = 2550.3333
$ws.Range("K13").Value = 2807.5713
$ws.Range("L13").Value = 7650.999899999999
$ws.Range("M13").Value = -2639.5713
$ws.Range("N13").Value = -7986.999899999999

$ws.Range("H38").Value = 83650
$ws.Range("I38").Value = 74
$ws.Range("J38").Value = 100365.2
$ws.Range("K38").Value = 222
$ws.Range("L38").Value = 301095.6
$ws.Range("M38").Value = 125
$ws.Range("N38").Value = -301789.6

$ws.Range("H113").Value = 1031.2916
$ws.Range("I113").Value = 1134.8334
$ws.Range("J113").Value = 720.6667
$ws.Range("K113").Value = 3404.5002
$ws.Range("L113").Value = 2162.0001
$ws.Range("M113").Value = -1234.5002
$ws.Range("N113").Value = -6502.0001

$ws.Range("H131").Value = 2376.7812
$ws.Range("I131").Value = 546
$ws.Range("J131").Value = 2477.3735
$ws.Range("K131").Value = 1638
$ws.Range("L131").Value = 7432.120500000001
$ws.Range("M131").Value = 3402
$ws.Range("N131").Value = -17512.1205

$ws.Range("H140").Value = 2189.1626
$ws.Range("I140").Value = 1801.9269
$ws.Range("J140").Value = 2596.2563
$ws.Range("K140").Value = 5405.780699999999
$ws.Range("L140").Value = 7788.7689
$ws.Range("M140").Value = -225.7806999999993
$ws.Range("N140").Value = -18148.7689

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2450
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 2900
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 2900
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -4896

$ws.Range("H83").Value = 2450
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 2900
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 14500
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -24484

$ws.Range("H132").Value = 2547.318
$ws.Range("I132").Value = 2244.1936
$ws.Range("J132").Value = 3270.1538
$ws.Range("K132").Value = 6732.5808
$ws.Range("L132").Value = 9810.4614
$ws.Range("M132").Value = -4202.5808
$ws.Range("N132").Value = -14870.4614

$ws.Range("H138").Value = 48884.93
$ws.Range("J138").Value = 48884.93
$ws.Range("L138").Value = 48884.93
$ws.Range("N138").Value = -59164.93

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1324.5625
$ws.Range("I22").Value = 888.2222
$ws.Range("J22").Value = 1885.5714
$ws.Range("K22").Value = 888.2222
$ws.Range("L22").Value = 1885.5714
$ws.Range("M22").Value = -593.2222
$ws.Range("N22").Value = -2475.5714

$ws.Range("H27").Value = 1324.5625
$ws.Range("I27").Value = 888.2222
$ws.Range("J27").Value = 1885.5714
$ws.Range("K27").Value = 888.2222
$ws.Range("L27").Value = 1885.5714
$ws.Range("M27").Value = -781.2222
$ws.Range("N27").Value = -2099.5714

$ws.Range("H46").Value = 1849.5
$ws.Range("I46").Value = 1710.4445
$ws.Range("J46").Value = 2266.6667
$ws.Range("K46").Value = 1710.4445
$ws.Range("L46").Value = 2266.6667
$ws.Range("M46").Value = -1522.4445
$ws.Range("N46").Value = -2642.6667

$ws.Range("H55").Value = 524
$ws.Range("I55").Value = 301
$ws.Range("J55").Value = 598.3333
$ws.Range("K55").Value = 301
$ws.Range("L55").Value = 598.3333
$ws.Range("M55").Value = -128
$ws.Range("N55").Value = -944.3333

$ws.Range("H82").Value = 2569
$ws.Range("I82").Value = 2372.5
$ws.Range("J82").Value = 2700
$ws.Range("K82").Value = 2372.5
$ws.Range("L82").Value = 2700
$ws.Range("M82").Value = -2011.5
$ws.Range("N82").Value = -3422

$ws.Range("H85").Value = 2569
$ws.Range("I85").Value = 2372.5
$ws.Range("J85").Value = 2700
$ws.Range("K85").Value = 2372.5
$ws.Range("L85").Value = 2700
$ws.Range("M85").Value = -1124.5
$ws.Range("N85").Value = -5196

$ws.Range("H127").Value = 53920.555
$ws.Range("J127").Value = 53920.555
$ws.Range("L127").Value = 53920.555
$ws.Range("N127").Value = -63840.555

$ws.Range("H132").Value = 7186.1177
$ws.Range("I132").Value = 2626.5
$ws.Range("J132").Value = 18129.2
$ws.Range("K132").Value = 7879.5
$ws.Range("L132").Value = 54387.60000000001
$ws.Range("M132").Value = -5349.5
$ws.Range("N132").Value = -59447.60000000001

$ws.Range("H133").Value = 24066
$ws.Range("J133").Value = 24066
$ws.Range("L133").Value = 24066
$ws.Range("N133").Value = -29126

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H132").Value = 1717.4329
$ws.Range("I132").Value = 1231.2
$ws.Range("J132").Value = 3147.5293
$ws.Range("K132").Value = 3693.6
$ws.Range("L132").Value = 9442.5879
$ws.Range("M132").Value = -1163.6
$ws.Range("N132").Value = -14502.5879
